# Actualización automática de grupos experimentales
# Updates the "Grupo_Experimental" (Con/Sin SmartScore) assignment for
# several participants and fixes the SmartScore numeric columns for
# participant "Emilio Rugerio" (row 12) which were stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap experimental group labels -----------------------------------
$ws.Range("B3").Value  = "Sin SmartScore"
$ws.Range("B7").Value  = "Con SmartScore"
$ws.Range("B10").Value = "Sin SmartScore"
$ws.Range("B11").Value = "Sin SmartScore"
$ws.Range("B12").Value = "Con SmartScore"

# --- Fix SmartScore values for row 12 (Emilio Rugerio) -----------------
# These were previously stored as text; convert them to real numbers.
$ws.Range("I12").Value  = 0.608
$ws.Range("L12").Value  = 0.58
$ws.Range("O12").Value  = 0.547
$ws.Range("R12").Value  = 0.688
$ws.Range("U12").Value  = 0.55
$ws.Range("X12").Value  = 0.515
$ws.Range("AA12").Value = 0.693
$ws.Range("AD12").Value = 0.474
$ws.Range("AG12").Value = 0.441
